$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new "F" column for "average debt per person" (it used to live in
# --- column E). Give F1 the same header formatting as E1 by copying it,
# --- then relabel E1 as the new "% of total people" column.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "average debt per person"
$ws.Range("E1").Value = "% of total people"

# --- Occupation rows, re-sorted (treasurer, doctors, merchant, farmer,
# --- executors) with the new "% of total people" column filled in
# --- (# of people / total # of people * 100) and "average debt per
# --- person" recomputed into its new column F.
# columns: B=occupation, C=6p_total, D=# of people, E=% of total people, F=average debt per person
$ws.Range("B2").Value = "treasurer"
$ws.Range("C2").Value = 12283.44
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 9.090909090909092
$ws.Range("F2").Value = 12283.44

$ws.Range("B3").Value = "doctors"
$ws.Range("C3").Value = 2008.37
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 9.090909090909092
$ws.Range("F3").Value = 2008.37

$ws.Range("B4").Value = "merchant"
$ws.Range("C4").Value = 861.53
$ws.Range("D4").Value = 7
$ws.Range("E4").Value = 63.63636363636363
$ws.Range("F4").Value = 123.0757142857143

$ws.Range("B5").Value = "farmer"
$ws.Range("C5").Value = 62.51
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 9.090909090909092
$ws.Range("F5").Value = 62.51

$ws.Range("B6").Value = "executors"
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 9.090909090909092
$ws.Range("F6").Value = 0
